# Update the "想去人数" (F column) values on the "展览" and "全部类型"
# worksheets to reflect the latest scrape numbers.

$wb = $excel.ActiveWorkbook

# Row -> new value map (same updates apply on both sheets).
$updates = @{
    4  = 83
    5  = 385
    6  = 11370
    7  = 692
    8  = 106
    12 = 159
    14 = 48
    15 = 46
    18 = 322
    19 = 1270
    21 = 892
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
